$d = $word.ActiveDocument

# Locate the paragraph that starts the real content ("A partir de los
# archivos..."). Everything before it -- the "Podrá encontrar más
# información..." sentence, the GitHub link, and the empty paragraph
# that only carried the _GoBack bookmark -- gets removed; that
# paragraph becomes the new first paragraph of the document.
$introEnd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "A partir de los archivos*") {
        $introEnd = $p.Range.Start
        break
    }
}
if ($introEnd -gt 0) {
    $d.Range(0, $introEnd).Delete()
}

# Fill in the word-count placeholder: "recopilar X palabras" ->
# "recopilar 32430 palabras".
$d.Content.Find.Execute("X", $true, $true, $false, $false, $false, $true, 1, $false, "32430", 2) | Out-Null

# Drop the accidental duplicated "de" (previously wrapped in proofErr
# spell-check tags) plus its extra space, joining back into a single
# "... tendencia central de sus longitudes:".
$d.Content.Find.Execute("tendencia central de de sus longitudes", $true, $false, $false, $false, $false, $true, 1, $false, "tendencia central de sus longitudes", 2) | Out-Null

# Re-anchor the _GoBack bookmark right before "sus longitudes:", where
# the duplicated word used to sit.
$target = $d.Content
$target.Find.Execute("sus longitudes:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $anchor) | Out-Null
